$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was bumped forward by one
# day (45204 -> 45205) for every data row (rows 2 through 61).
$ws.Range("C2:C61").Value = 45205
